$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.44917646495029
$ws.Range("C2").Value = 11.09444116004228
$ws.Range("D2").Value = 10.49533912789452
$ws.Range("F2").Value = 33.46927583699289
$ws.Range("G2").Value = 34.80690967482921
$ws.Range("H2").Value = 15.40028653559042
$ws.Range("I2").Value = 20.96967096590631
$ws.Range("J2").Value = 11.33186244100798
$ws.Range("N2").Value = 16.79162394413106
$ws.Range("B3").Value = 14.86898759129795
$ws.Range("C3").Value = 10.54071566256679
$ws.Range("D3").Value = 10.44192863122314
$ws.Range("F3").Value = 33.36177026009133
$ws.Range("G3").Value = 34.56357270112459
$ws.Range("H3").Value = 15.42843967714146
$ws.Range("I3").Value = 21.05269386136897
$ws.Range("J3").Value = 11.31532167394267
$ws.Range("N3").Value = 16.85279672142622
$ws.Range("B4").Value = 14.50365563295848
$ws.Range("C4").Value = 10.18720076475601
$ws.Range("D4").Value = 10.41081180161269
$ws.Range("F4").Value = 33.30660983154578
$ws.Range("G4").Value = 34.42812223927739
$ws.Range("H4").Value = 15.44955248768197
$ws.Range("I4").Value = 21.10955180328134
$ws.Range("J4").Value = 11.30757151028261
$ws.Range("N4").Value = 16.89224534767201
$ws.Range("B5").Value = 14.35272888074675
$ws.Range("C5").Value = 10.03992388208138
$ws.Range("D5").Value = 10.39856233888271
$ws.Range("F5").Value = 33.2868712435663
$ws.Range("G5").Value = 34.37648442098389
$ws.Range("H5").Value = 15.45911502735422
$ws.Range("I5").Value = 21.13419314976219
$ws.Range("J5").Value = 11.3050199003665
$ws.Range("N5").Value = 16.90879714659867
$ws.Range("B6").Value = 14.32755109029291
$ws.Range("C6").Value = 10.015280371576
$ws.Range("D6").Value = 10.39655460456612
$ws.Range("F6").Value = 33.28375940907267
$ws.Range("G6").Value = 34.36812625280085
$ws.Range("H6").Value = 15.46076069129634
$ws.Range("I6").Value = 21.13837346052379
$ws.Range("J6").Value = 11.30463288651049
$ws.Range("N6").Value = 16.91157435886655
$ws.Range("B7").Value = 14.5016281574725
$ws.Range("C7").Value = 10.18522729568199
$ws.Range("D7").Value = 10.410644844802
$ws.Range("F7").Value = 33.30633252433554
$ws.Range("G7").Value = 34.42741136318752
$ws.Range("H7").Value = 15.4496775735956
$ws.Range("I7").Value = 21.10987817830744
$ws.Range("J7").Value = 11.30753464026323
$ws.Range("N7").Value = 16.89246664129803
$ws.Range("B8").Value = 15.25114872329513
$ws.Range("C8").Value = 10.9064252998534
$ws.Range("D8").Value = 10.4765806625856
$ws.Range("F8").Value = 33.42996456057222
$ws.Range("G8").Value = 34.72014174943762
$ws.Range("H8").Value = 15.40919758407444
$ws.Range("I8").Value = 20.99707196028341
$ws.Range("J8").Value = 11.32566095964947
$ws.Range("N8").Value = 16.8123252057701
$ws.Range("B9").Value = 16.63964913505084
$ws.Range("C9").Value = 12.20669493368631
$ws.Range("D9").Value = 10.61876607015449
$ws.Range("F9").Value = 33.75778806123581
$ws.Range("G9").Value = 35.40239990945884
$ws.Range("H9").Value = 15.36031747194778
$ws.Range("I9").Value = 20.82286510266748
$ws.Range("J9").Value = 11.38020503935202
$ws.Range("N9").Value = 16.67008991651754
$ws.Range("B10").Value = 17.59955552870298
$ws.Range("C10").Value = 13.0852597636549
$ws.Range("D10").Value = 10.73048379638974
$ws.Range("F10").Value = 34.04950749177087
$ws.Range("G10").Value = 35.9656996156148
$ws.Range("H10").Value = 15.34317088029587
$ws.Range("I10").Value = 20.72398550976802
$ws.Range("J10").Value = 11.43171029780734
$ws.Range("N10").Value = 16.57460118294977
$ws.Range("B11").Value = 18.02133626073326
$ws.Range("C11").Value = 13.46715497217175
$ws.Range("D11").Value = 10.78274476030848
$ws.Range("F11").Value = 34.19293520723063
$ws.Range("G11").Value = 36.23446098255506
$ws.Range("H11").Value = 15.33947513357533
$ws.Range("I11").Value = 20.68541749747467
$ws.Range("J11").Value = 11.45758201836099
$ws.Range("N11").Value = 16.53309980864204
$ws.Range("B12").Value = 18.17878623648515
$ws.Range("C12").Value = 13.60914165236536
$ws.Range("D12").Value = 10.80272949303357
$ws.Range("F12").Value = 34.24875650272421
$ws.Range("G12").Value = 36.33794411099772
$ws.Range("H12").Value = 15.33866767490058
$ws.Range("I12").Value = 20.67174236785828
$ws.Range("J12").Value = 11.46772577203788
$ws.Range("N12").Value = 16.5176615462932
$ws.Range("B13").Value = 18.14497936086687
$ws.Range("C13").Value = 13.5786802244958
$ws.Range("D13").Value = 10.79841696788703
$ws.Range("F13").Value = 34.23666787392199
$ws.Range("G13").Value = 36.31558264242308
$ws.Range("H13").Value = 15.33881521906951
$ws.Range("I13").Value = 20.6746460829335
$ws.Range("J13").Value = 11.46552578905221
$ws.Range("N13").Value = 16.52097413210582
$ws.Range("B14").Value = 18.03433592056825
$ws.Range("C14").Value = 13.47888935970704
$ws.Range("D14").Value = 10.78438508775234
$ws.Range("F14").Value = 34.19749763360013
$ws.Range("G14").Value = 36.24294088320685
$ws.Range("H14").Value = 15.33939682887869
$ws.Range("I14").Value = 20.68427376033794
$ws.Range("J14").Value = 11.45840962973411
$ws.Range("N14").Value = 16.53182414017747
$ws.Range("B15").Value = 17.96626456903221
$ws.Range("C15").Value = 13.41742027896901
$ws.Range("D15").Value = 10.77581512948782
$ws.Range("F15").Value = 34.17370010682227
$ws.Range("G15").Value = 36.19866551834711
$ws.Range("H15").Value = 15.33983023108446
$ws.Range("I15").Value = 20.69029229148579
$ws.Range("J15").Value = 11.45409578760914
$ws.Range("N15").Value = 16.53850617722368
$ws.Range("B16").Value = 17.57168023362612
$ws.Range("C16").Value = 13.05993785361117
$ws.Range("D16").Value = 10.72709634091094
$ws.Range("F16").Value = 34.040347127968
$ws.Range("G16").Value = 35.94837960559764
$ws.Range("H16").Value = 15.34349516170536
$ws.Range("I16").Value = 20.72663571578485
$ws.Range("J16").Value = 11.43006828605172
$ws.Range("N16").Value = 16.5773522749353
$ws.Range("B17").Value = 17.32570586258715
$ws.Range("C17").Value = 12.83602616008632
$ws.Range("D17").Value = 10.69756902351349
$ws.Range("F17").Value = 33.96126293983744
$ws.Range("G17").Value = 35.79797943124156
$ws.Range("H17").Value = 15.34679612439137
$ws.Range("I17").Value = 20.75057917192565
$ws.Range("J17").Value = 11.41595057368753
$ws.Range("N17").Value = 16.60167838600277
$ws.Range("B18").Value = 17.18283335229236
$ws.Range("C18").Value = 12.70556907698046
$ws.Range("D18").Value = 10.68072199852665
$ws.Range("F18").Value = 33.91678734054906
$ws.Range("G18").Value = 35.71265678017464
$ws.Range("H18").Value = 15.34908098260453
$ws.Range("I18").Value = 20.76495397894096
$ws.Range("J18").Value = 11.40806055914245
$ws.Range("N18").Value = 16.61585250007029
$ws.Range("B19").Value = 17.13422386532058
$ws.Range("C19").Value = 12.66111437238957
$ws.Range("D19").Value = 10.67504166958213
$ws.Range("F19").Value = 33.90190337014419
$ws.Range("G19").Value = 35.68397395355055
$ws.Range("H19").Value = 15.34992086963997
$ws.Range("I19").Value = 20.76992440077238
$ws.Range("J19").Value = 11.4054287864172
$ws.Range("N19").Value = 16.62068296975884
$ws.Range("B20").Value = 17.35203557282473
$ws.Range("C20").Value = 12.8600352545583
$ws.Range("D20").Value = 10.7006982386399
$ws.Range("F20").Value = 33.96957713037964
$ws.Range("G20").Value = 35.8138679346667
$ws.Range("H20").Value = 15.3464047433411
$ws.Range("I20").Value = 20.74796786945958
$ws.Range("J20").Value = 11.41742964231616
$ws.Range("N20").Value = 16.59906996275868
$ws.Range("B21").Value = 18.06689710874684
$ws.Range("C21").Value = 13.50827221915199
$ws.Range("D21").Value = 10.78850140925273
$ws.Range("F21").Value = 34.20896223086592
$ws.Range("G21").Value = 36.26423188765903
$ws.Range("H21").Value = 15.33920991534963
$ws.Range("I21").Value = 20.68142058592364
$ws.Range("J21").Value = 11.46049044650701
$ws.Range("N21").Value = 16.52862970984124
$ws.Range("B22").Value = 18.52082559302202
$ws.Range("C22").Value = 13.91658516420674
$ws.Range("D22").Value = 10.84701438983019
$ws.Range("F22").Value = 34.37418813912803
$ws.Range("G22").Value = 36.56848755004633
$ws.Range("H22").Value = 15.33795913058395
$ws.Range("I22").Value = 20.64335067609168
$ws.Range("J22").Value = 11.49065160670293
$ws.Range("N22").Value = 16.48420937376746
$ws.Range("B23").Value = 18.27980826439908
$ws.Range("C23").Value = 13.70008607892276
$ws.Range("D23").Value = 10.81568583590726
$ws.Range("F23").Value = 34.28521307603128
$ws.Range("G23").Value = 36.40522411583238
$ws.Range("H23").Value = 15.33831037897184
$ws.Range("I23").Value = 20.66317069072462
$ws.Range("J23").Value = 11.4743709174289
$ws.Range("N23").Value = 16.50776980217178
$ws.Range("B24").Value = 17.34013644397404
$ws.Range("C24").Value = 12.84918611730963
$ws.Range("D24").Value = 10.69928311890327
$ws.Range("F24").Value = 33.96581519559069
$ws.Range("G24").Value = 35.80668117072827
$ws.Range("H24").Value = 15.3465804810888
$ws.Range("I24").Value = 20.74914654165425
$ws.Range("J24").Value = 11.41676024971953
$ws.Range("N24").Value = 16.60024864299894
$ws.Range("B25").Value = 16.2739124316183
$ws.Range("C25").Value = 11.86797406010182
$ws.Range("D25").Value = 10.57898182768162
$ws.Range("F25").Value = 33.66006991399135
$ws.Range("G25").Value = 35.2066245247873
$ws.Range("H25").Value = 15.37025648616859
$ws.Range("I25").Value = 20.86491282420375
$ws.Range("J25").Value = 11.36342913597897
$ws.Range("N25").Value = 16.70697975068731
